$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new rows before row 270, pushing existing rows 270-274 down to 272-276
$ws.Rows.Item(270).Resize(2).Insert()

# New row 270
$ws.Cells.Item(270, 1).Value = 3
$ws.Cells.Item(270, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(270, 3).Value = "Coquimbo"
$ws.Cells.Item(270, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(270, 4).Value = 44448
$ws.Cells.Item(270, 5).Value = 5
$ws.Cells.Item(270, 6).Value = 100114001
$ws.Cells.Item(270, 7).Value = "Papa"
$ws.Cells.Item(270, 8).Value = "Asterix"
$ws.Cells.Item(270, 9).Value = "1a (guarda)"
$ws.Cells.Item(270, 10).Value = 450
$ws.Cells.Item(270, 11).Value = 8500
$ws.Cells.Item(270, 12).Value = 9000
$ws.Cells.Item(270, 13).Value = 8789
$ws.Cells.Item(270, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(270, 15).Value = "Provincia de Talca"
$ws.Cells.Item(270, 16).Value = 352
$ws.Cells.Item(270, 17).Value = 25
$ws.Cells.Item(270, 18).Value = "Hortaliza"

# New row 271
$ws.Cells.Item(271, 1).Value = 3
$ws.Cells.Item(271, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(271, 3).Value = "Coquimbo"
$ws.Cells.Item(271, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(271, 4).Value = 44448
$ws.Cells.Item(271, 5).Value = 5
$ws.Cells.Item(271, 6).Value = 100114001
$ws.Cells.Item(271, 7).Value = "Papa"
$ws.Cells.Item(271, 8).Value = "Rosara"
$ws.Cells.Item(271, 9).Value = "1a (guarda)"
$ws.Cells.Item(271, 10).Value = 510
$ws.Cells.Item(271, 11).Value = 7500
$ws.Cells.Item(271, 12).Value = 8000
$ws.Cells.Item(271, 13).Value = 7725
$ws.Cells.Item(271, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(271, 15).Value = "Provincia de Talca"
$ws.Cells.Item(271, 16).Value = 309
$ws.Cells.Item(271, 17).Value = 25
$ws.Cells.Item(271, 18).Value = "Hortaliza"
